# Add a "credit" column (E) and replace the sid/password sample values
# with the new "0"-prefixed ones, matching the commit:
# "modify on recieving student credit inputs".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for the credit column.
$ws.Range("E1").Value = "credit"

# Update passwords first (Cathy, then Dorthy)...
$ws.Range("C2").Value = "spw033333"
$ws.Range("C3").Value = "spw044444"

# ...then the student ids (Cathy, then Dorthy).
$ws.Range("B2").Value = "sid033333"
$ws.Range("B3").Value = "sid044444"

# Record a credit value of 3 for each student.
$ws.Range("E2").Value = 3
$ws.Range("E3").Value = 3

# Leave the selection where the user's last entry was.
$ws.Range("B3").Select()
